# Apply the "Add cantrals by cantons" edit:
#  - The worksheet had two header rows (row 1 and row 2). They are replaced
#    by a single header row (row 1) with 11 columns (idx, idx2, Name,
#    Date Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter,
#    (GWh) Summer, (GWh) Year); all data rows shift up by one row.
#  - The view selection/scroll position is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Delete the old second header row (row 2). This shifts every data row
#    that followed it up by one position, so former row 3 (first data row)
#    becomes row 2, former row 51 (last data row) becomes row 50, etc.
$ws.Rows(2).Delete()

# The worksheet originally had a trailing pre-formatted blank row (row 136,
# columns L:O only) below the last real data row. Row-deletion above pulled
# everything up by one, so that trailing blank formatting now sits on row
# 135 and row 136 is plain again. Restore the same trailing formatted blank
# row at 136 so the sheet keeps its original 136-row footprint.
$ws.Range("L135:O135").Copy() | Out-Null
$ws.Range("L136:O136").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Rewrite row 1 as the single, unified header row.
$ws.Cells.Item(1,1).Value2 = "idx"
$ws.Cells.Item(1,2).Value2 = "idx2"
$ws.Cells.Item(1,3).Value2 = "Name"
$ws.Cells.Item(1,4).Value2 = "Date Start"
$ws.Cells.Item(1,5).Value2 = "Date End"
$ws.Cells.Item(1,6).Value2 = "(m3/s)"
$ws.Cells.Item(1,7).Value2 = "(MW1)"
$ws.Cells.Item(1,8).Value2 = "(MW2)"
$ws.Cells.Item(1,9).Value2 = "(GWh) Winter"
$ws.Cells.Item(1,10).Value2 = "(GWh) Summer"
$ws.Cells.Item(1,11).Value2 = "(GWh) Year"

# Columns A:E are plain (no special formatting, like the rest of row 1 used
# to be), while E1 previously carried the old "Année" header style, so make
# sure it is reset back to the default look.
$ws.Range("A1:E1").Style = "Normal"

# Columns F:K keep the same look as the rest of the header text (Arial 9).
for ($c = 6; $c -le 11; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 9
}

# 3) Update the saved view: the first data row (now row 2) is selected
#    instead of the old row 45.
$ws.Range("A2:K2").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
